$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (new): 21 loka - add the date first
$ws.Rows.Item(12).RowHeight = 29
$ws.Cells.Item(12, 1).Value = "21 loka"

# Row 10 (19 loka): shorten the "camera class" note - drop "ja ruudukko"
$ws.Cells.Item(10, 3).Value = "Kameraluokka, liikkuminen scenessä"

# Row 11 (20 loka): add a note about the infinite grid tutorial
$ws.Cells.Item(11, 3).Value = "Infinite grid tutoriaalia"
$ws.Cells.Item(11, 3).WrapText = $true

# Row 12 (new): fill in the hours worked
$ws.Cells.Item(12, 2).Value = "9.15-10.15, 14.00-17.00, "
$ws.Cells.Item(12, 2).NumberFormat = "h:mm"
$ws.Cells.Item(12, 2).WrapText = $true

$ws.Range("B12").Select()
